# feat: change PFAS selection criteria
#
# The PFAS substances considered for the sumPFAS aggregate changed, so
# "PFPeA" and "6:2 FTSA" are dropped from the per-(type, grp) breakdown
# table, and the resulting "sumPFAS" rows are recomputed against the new
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows for the PFAS substances no longer included
# (delete from the bottom up so earlier row numbers stay valid).
$rowsToDelete = @(56, 45, 42, 31, 28, 17, 14, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Recompute the sumPFAS rows (min, median, max) for the new selection.
$ws.Range("D2").Value = 3.280028437291977
$ws.Range("E2").Value = 4.836746552908711
$ws.Range("F2").Value = 8.785980259008392

$ws.Range("D14").Value = 10.69454611693446
$ws.Range("E14").Value = 14.67051176672183
$ws.Range("F14").Value = 17.48602230973641

$ws.Range("D26").Value = 13.99735378205541
$ws.Range("E26").Value = 20.81744502954255
$ws.Range("F26").Value = 25.4647628633788

$ws.Range("D38").Value = 14.21090609082056
$ws.Range("E38").Value = 17.0321208667622
$ws.Range("F38").Value = 25.71497209101526
